$d = $word.ActiveDocument

# 1. Remove the standalone "Meta description" paragraph near the top of the document
#    (it contained the bold "Meta description" label plus the description text).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new paragraph right before the final paragraph of the document,
#    containing a bold copy of the page title ("Play Big Bang ... | Review").
#    Build it from raw OOXML so the run/paragraph structure comes out clean
#    (no inherited formatting from neighboring paragraphs).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Big Bang (Belatra Games) Slot Game for Free | Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titlePara = $d.Paragraphs.Item($count)
$titlePara.Range.InsertXML($titleXml) | Out-Null

# 3. Replace the original final (italic) paragraph's text -- the old AI-art prompt --
#    with the meta-description copy that used to live at the top of the document.
$d.Content.Find.Execute(
    "Create a fun and colorful cartoon image featuring a happy Maya warrior wearing glasses as the main character. The background of the image should depict a planetary system with colorful planets and stars in the distance. The warrior should be holding a futuristic-looking device and surrounded by symbols from the game Big Bang (such as planets and stars). The overall tone of the image should be exciting and adventurous, inviting players to join the Maya warrior on a journey through space in the game. The image should be bright and eye-catching, with bold colors that pop off the screen.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Embark on a journey through the cosmos with Big Bang slot game by Belatra Games. Read our review and play for free.",
    2
) | Out-Null
